$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename measure labels in column A (mean_Beta_bb -> mean_frac, var_Beta_bb -> var_frac)
for ($i = 2; $i -le 155; $i++) {
    $cell = $ws.Cells.Item($i, 1)
    $v = $cell.Value2
    if ($v -eq "mean_Beta_bb") {
        $cell.Value = "mean_frac"
    }
    if ($v -eq "var_Beta_bb") {
        $cell.Value = "var_frac"
    }
}

# Add new columns L (vTest) and M (vstat)
$ws.Range("L1").Value = "vTest"
$ws.Range("M1").Value = "vstat"

$ws.Range("L2").Value = [double]"2.5762827336064148e-05"
$ws.Range("M2").Value = [double]"5.211860363686708"
$ws.Range("L3").Value = [double]"0.60797285644484833"
$ws.Range("M3").Value = [double]"0.26597362257314955"
$ws.Range("L4").Value = [double]"9.6613304174484718e-05"
$ws.Range("M4").Value = [double]"4.5659344681198792"
$ws.Range("L5").Value = [double]"0.011011177873634358"
$ws.Range("M5").Value = [double]"6.891675377058367"
$ws.Range("L6").Value = [double]"7.3971141519727733e-08"
$ws.Range("M6").Value = [double]"37.771363245396259"
$ws.Range("L7").Value = [double]"0.73660592620094012"
$ws.Range("M7").Value = [double]"0.11420524379701054"
$ws.Range("L8").Value = [double]"2.576230913399782e-06"
$ws.Range("M8").Value = [double]"27.086338579426272"
$ws.Range("L9").Value = [double]"0.43742099062190409"
$ws.Range("M9").Value = [double]"0.61131139688638658"
$ws.Range("L10").Value = [double]"0.0067690336629330355"
$ws.Range("M10").Value = [double]"7.878429023907489"
$ws.Range("L11").Value = [double]"0.5985994644217757"
$ws.Range("M11").Value = [double]"0.28013369169253893"
$ws.Range("L12").Value = [double]"0.0049605777876470464"
$ws.Range("M12").Value = [double]"8.5221394450079924"
$ws.Range("L13").Value = [double]"0.53017873857353448"
$ws.Range("M13").Value = [double]"0.39874185498002757"
$ws.Range("L14").Value = [double]"1.4075935121890263e-06"
$ws.Range("M14").Value = [double]"28.811531435330643"
$ws.Range("L15").Value = [double]"3.3856144825795643e-08"
$ws.Range("M15").Value = [double]"40.308415308426675"
$ws.Range("L16").Value = [double]"2.660617396146009e-06"
$ws.Range("M16").Value = [double]"26.995394321103845"
$ws.Range("L17").Value = [double]"1.2639919809689778e-08"
$ws.Range("M17").Value = [double]"43.606907213582566"
$ws.Range("L18").Value = [double]"6.033301949308423e-05"
$ws.Range("M18").Value = [double]"18.676548487400254"
$ws.Range("L19").Value = [double]"1.9037964363741201e-05"
$ws.Range("M19").Value = [double]"21.639593561217623"
$ws.Range("L20").Value = [double]"2.9804822730551707e-08"
$ws.Range("M20").Value = [double]"40.728748370983333"
$ws.Range("L21").Value = [double]"0.0050532790873714225"
$ws.Range("M21").Value = [double]"8.483512651930397"
$ws.Range("L22").Value = [double]"0.52475999502224147"
$ws.Range("M22").Value = [double]"0.40938240793386121"
$ws.Range("L23").Value = [double]"0.00087381814682843727"
$ws.Range("M23").Value = [double]"12.29768698448647"
$ws.Range("L24").Value = [double]"8.2866848301519554e-05"
$ws.Range("M24").Value = [double]"17.883697723718818"
$ws.Range("L25").Value = [double]"0.10680206774112666"
$ws.Range("M25").Value = [double]"2.6821441289466019"
$ws.Range("L26").Value = [double]"0.015651550953533776"
$ws.Range("M26").Value = [double]"0.4026940229975135"
$ws.Range("L27").Value = [double]"6.3707260594425491e-07"
$ws.Range("M27").Value = [double]"31.131447908076005"
$ws.Range("L28").Value = [double]"0.51788710866703336"
$ws.Range("M28").Value = [double]"0.42316531686009668"
$ws.Range("L29").Value = [double]"0.25445086818006002"
$ws.Range("M29").Value = [double]"1.3243766212756993"
$ws.Range("L30").Value = [double]"0.35458586500134481"
$ws.Range("M30").Value = [double]"0.87062273747392427"
$ws.Range("L31").Value = [double]"0.28875930919129861"
$ws.Range("M31").Value = [double]"0.67431340717246246"
$ws.Range("L32").Value = [double]"0.51857668651447941"
$ws.Range("M32").Value = [double]"0.42176776082392098"
$ws.Range("L33").Value = [double]"0.38353981526739744"
$ws.Range("M33").Value = [double]"0.77077123142753989"
$ws.Range("L34").Value = [double]"0.31612897087781539"
$ws.Range("M34").Value = [double]"0.68920353971527826"
$ws.Range("L35").Value = [double]"0.28278506429731232"
$ws.Range("M35").Value = [double]"1.1750050331482469"
$ws.Range("L36").Value = [double]"0.10749475111401197"
$ws.Range("M36").Value = [double]"0.54857425926311065"
$ws.Range("L37").Value = [double]"0.35231513039407159"
$ws.Range("M37").Value = [double]"0.87893656518932572"
$ws.Range("L38").Value = [double]"0.23271607672041322"
$ws.Range("M38").Value = [double]"1.5626609901801569"
$ws.Range("L39").Value = [double]"0.22477006507260022"
$ws.Range("M39").Value = [double]"1.5050442812333489"
$ws.Range("L40").Value = [double]"0.0010242218650618183"
$ws.Range("M40").Value = [double]"11.940092341542162"
$ws.Range("L41").Value = [double]"0.51125745814367096"
$ws.Range("M41").Value = [double]"0.43677114826960411"
$ws.Range("L42").Value = [double]"0.68063049349558069"
$ws.Range("M42").Value = [double]"0.17110601089592067"
$ws.Range("L43").Value = [double]"0.011840715358962951"
$ws.Range("M43").Value = [double]"6.7466228097155856"
$ws.Range("L44").Value = [double]"0.13847244326686414"
$ws.Range("M44").Value = [double]"2.2555171670420493"
$ws.Range("L45").Value = [double]"0.053944581379353387"
$ws.Range("M45").Value = [double]"0.48618563457995112"
$ws.Range("L46").Value = [double]"0.063531593923763396"
$ws.Range("M46").Value = [double]"2.009889364748557"
$ws.Range("L47").Value = [double]"0.26359814783146035"
$ws.Range("M47").Value = [double]"1.5188767286542639"
$ws.Range("L48").Value = [double]"0.043766858973133337"
$ws.Range("M48").Value = [double]"4.2457395125359456"
$ws.Range("L49").Value = [double]"0.033217572398606718"
$ws.Range("M49").Value = [double]"0.44996930770268268"
$ws.Range("L50").Value = [double]"0.22603831110560235"
$ws.Range("M50").Value = [double]"1.4967343165996188"
$ws.Range("L51").Value = [double]"0.091732352566521222"
$ws.Range("M51").Value = [double]"1.8838746049957995"
$ws.Range("L52").Value = [double]"0.19981622148749531"
$ws.Range("M52").Value = [double]"0.6206484512456989"
$ws.Range("L53").Value = [double]"0.027013781177575441"
$ws.Range("M53").Value = [double]"0.43601741098586089"
$ws.Range("L54").Value = [double]"0.084761876984754361"
$ws.Range("M54").Value = [double]"3.0737096675529862"
$ws.Range("L55").Value = [double]"0.081252485946917724"
$ws.Range("M55").Value = [double]"3.1464550112963892"
$ws.Range("L56").Value = [double]"0.048626982658208043"
$ws.Range("M56").Value = [double]"4.0543184020705416"
$ws.Range("L57").Value = [double]"0.24227424788197027"
$ws.Range("M57").Value = [double]"1.3951652955623326"
$ws.Range("L58").Value = [double]"0.98947723808961063"
$ws.Range("M58").Value = [double]"0.00017542257267018368"
$ws.Range("L59").Value = [double]"0.0071048525889218081"
$ws.Range("M59").Value = [double]"0.36172574499984539"
$ws.Range("L60").Value = [double]"0.57922962130253908"
$ws.Range("M60").Value = [double]"0.31091364454406933"
$ws.Range("L61").Value = [double]"0.10810928726810623"
$ws.Range("M61").Value = [double]"1.8275573869079906"
$ws.Range("L62").Value = [double]"0.79276728641820005"
$ws.Range("M62").Value = [double]"1.1034058931377237"
$ws.Range("L63").Value = [double]"0.28040038340483908"
$ws.Range("M63").Value = [double]"1.186826466513152"
$ws.Range("L64").Value = [double]"0.0049927450323955591"
$ws.Range("M64").Value = [double]"8.508650762862203"
$ws.Range("L65").Value = [double]"0.61951442826479819"
$ws.Range("M65").Value = [double]"0.24917275531847358"
$ws.Range("L66").Value = [double]"0.32392710315357276"
$ws.Range("M66").Value = [double]"1.4456416311997167"
$ws.Range("L67").Value = [double]"0.045284105962053282"
$ws.Range("M67").Value = [double]"4.1836071540067561"
$ws.Range("L68").Value = [double]"0.17692170516270189"
$ws.Range("M68").Value = [double]"1.8677227208475142"
$ws.Range("L69").Value = [double]"0.99090743166539552"
$ws.Range("M69").Value = [double]"1.0050219091703165"
$ws.Range("L70").Value = [double]"0.20176176297531856"
$ws.Range("M70").Value = [double]"1.6665130843437674"
$ws.Range("L71").Value = [double]"0.16482321304927194"
$ws.Range("M71").Value = [double]"1.9782230050853971"
$ws.Range("L72").Value = [double]"0.37077299242454276"
$ws.Range("M72").Value = [double]"0.71741926284954816"
$ws.Range("L73").Value = [double]"0.2096462036716531"
$ws.Range("M73").Value = [double]"1.608774061292771"
$ws.Range("L74").Value = [double]"0.83604922835154039"
$ws.Range("M74").Value = [double]"0.043208398892365746"
$ws.Range("L75").Value = [double]"0.4776573707082793"
$ws.Range("M75").Value = [double]"0.76863568348141931"
$ws.Range("L76").Value = [double]"0.63289883475360842"
$ws.Range("M76").Value = [double]"0.8380672625941209"
$ws.Range("L77").Value = [double]"0.36487024873779211"
$ws.Range("M77").Value = [double]"0.71445348666315778"
$ws.Range("L78").Value = [double]"0.34794664502136163"
$ws.Range("M78").Value = [double]"1.4199291352950729"
$ws.Range("L79").Value = [double]"0.36592300322208815"
$ws.Range("M79").Value = [double]"0.71498374149112021"
$ws.Range("L80").Value = [double]"0.3319601885567342"
$ws.Range("M80").Value = [double]"0.6975661275818692"
$ws.Range("L81").Value = [double]"0.5855632711456279"
$ws.Range("M81").Value = [double]"1.2257973770674342"
$ws.Range("L82").Value = [double]"0.31275911122967931"
$ws.Range("M82").Value = [double]"0.68740106445525118"
$ws.Range("L83").Value = [double]"0.03236317275843853"
$ws.Range("M83").Value = [double]"4.80339327756076"
$ws.Range("L84").Value = [double]"0.3587151510551389"
$ws.Range("M84").Value = [double]"0.71134160746089503"
$ws.Range("L85").Value = [double]"0.11325337929100232"
$ws.Range("M85").Value = [double]"2.5845004104041873"
$ws.Range("L86").Value = [double]"0.33985797140851865"
$ws.Range("M86").Value = [double]"1.4284031450438801"
$ws.Range("L87").Value = [double]"0.0051271108253447099"
$ws.Range("M87").Value = [double]"2.9062690608690098"
$ws.Range("L88").Value = [double]"0.47070820453966755"
$ws.Range("M88").Value = [double]"0.52707801591151826"
$ws.Range("L89").Value = [double]"0.85703850768011203"
$ws.Range("M89").Value = [double]"0.032737567001683454"
$ws.Range("L90").Value = [double]"0.23079877118144487"
$ws.Range("M90").Value = [double]"1.4660447862359247"
$ws.Range("L91").Value = [double]"0.70622132790525449"
$ws.Range("M91").Value = [double]"0.14346360815353487"
$ws.Range("L92").Value = [double]"0.48426524823317429"
$ws.Range("M92").Value = [double]"0.49547232388038687"
$ws.Range("L93").Value = [double]"0.066596713707228708"
$ws.Range("M93").Value = [double]"0.50370998376610909"
$ws.Range("L94").Value = [double]"0.44997843544955218"
$ws.Range("M94").Value = [double]"1.3258414850544891"
$ws.Range("L95").Value = [double]"0.05125769463202233"
$ws.Range("M95").Value = [double]"3.9591736208292772"
$ws.Range("L96").Value = [double]"0.78395880291417042"
$ws.Range("M96").Value = [double]"0.075854072789653296"
$ws.Range("L97").Value = [double]"0.16195400347406638"
$ws.Range("M97").Value = [double]"1.6885889912917023"
$ws.Range("L98").Value = [double]"0.056139991723807375"
$ws.Range("M98").Value = [double]"3.7958912023477418"
$ws.Range("L99").Value = [double]"0.0047986186277598286"
$ws.Range("M99").Value = [double]"8.5914791395320442"
$ws.Range("L100").Value = [double]"0.36265043913859368"
$ws.Range("M100").Value = [double]"1.4049696276060155"
$ws.Range("L101").Value = [double]"0.012741747485159006"
$ws.Range("M101").Value = [double]"0.39128209116628465"
$ws.Range("L102").Value = [double]"0.74706965338239984"
$ws.Range("M102").Value = [double]"1.1282161136306774"
$ws.Range("L103").Value = [double]"0.49390110964889028"
$ws.Range("M103").Value = [double]"0.7761021721597513"
$ws.Range("L104").Value = [double]"0.13300926645080036"
$ws.Range("M104").Value = [double]"2.3206478174043608"
$ws.Range("L105").Value = [double]"0.88017187251673956"
$ws.Range("M105").Value = [double]"0.022923814440688459"
$ws.Range("L106").Value = [double]"0.03175506529895588"
$ws.Range("M106").Value = [double]"4.8388624023373428"
$ws.Range("L107").Value = [double]"0.73378589029522268"
$ws.Range("M107").Value = [double]"1.1356225741682915"
$ws.Range("L108").Value = [double]"5.0805215591612466e-05"
$ws.Range("M108").Value = [double]"0.20834318271086447"
$ws.Range("L109").Value = [double]"0.0049311205302741161"
$ws.Range("M109").Value = [double]"0.34509883907885325"
$ws.Range("L110").Value = [double]"0.030465515453339375"
$ws.Range("M110").Value = [double]"0.444034435170299"
$ws.Range("L111").Value = [double]"0.004470332840056864"
$ws.Range("M111").Value = [double]"8.7398923541622473"
$ws.Range("L112").Value = [double]"0.3623388846407245"
$ws.Range("M112").Value = [double]"0.71317611457222219"
$ws.Range("L113").Value = [double]"0.026016834748084961"
$ws.Range("M113").Value = [double]"5.2144576988207314"
$ws.Range("L114").Value = [double]"0.25310802200005755"
$ws.Range("M114").Value = [double]"1.3319733771751774"
$ws.Range("L115").Value = [double]"0.64812102641679537"
$ws.Range("M115").Value = [double]"0.21042356142408886"
$ws.Range("L116").Value = [double]"0.057922470655001924"
$ws.Range("M116").Value = [double]"3.7400989895579584"
$ws.Range("L117").Value = [double]"0.20477407060293196"
$ws.Range("M117").Value = [double]"1.6441376559026058"
$ws.Range("L118").Value = [double]"0.67300331464306273"
$ws.Range("M118").Value = [double]"0.85557437081028143"
$ws.Range("L119").Value = [double]"0.28385548605553024"
$ws.Range("M119").Value = [double]"1.1697403651053524"
$ws.Range("L120").Value = [double]"0.0069383153001743277"
$ws.Range("M120").Value = [double]"7.8277150638174877"
$ws.Range("L121").Value = [double]"0.10594147341694665"
$ws.Range("M121").Value = [double]"0.54709180854606343"
$ws.Range("L122").Value = [double]"0.12414886188486585"
$ws.Range("M122").Value = [double]"2.4330643089372406"
$ws.Range("L123").Value = [double]"0.49962167005815006"
$ws.Range("M123").Value = [double]"0.77871666931627703"
$ws.Range("L124").Value = [double]"0.26598985471722292"
$ws.Range("M124").Value = [double]"1.2611112301087131"
$ws.Range("L125").Value = [double]"0.041027273571029331"
$ws.Range("M125").Value = [double]"4.3640597108529971"
$ws.Range("L126").Value = [double]"0.11032972439689885"
$ws.Range("M126").Value = [double]"0.55124621607573132"
$ws.Range("L127").Value = [double]"0.22501409514209903"
$ws.Range("M127").Value = [double]"1.5034408629805209"
$ws.Range("L128").Value = [double]"0.52940542428588444"
$ws.Range("M128").Value = [double]"0.79221645634542381"
$ws.Range("L129").Value = [double]"0.17042075321645367"
$ws.Range("M129").Value = [double]"1.9259545160802316"
$ws.Range("L130").Value = [double]"0.66083633447071"
$ws.Range("M130").Value = [double]"0.19446687932616091"
$ws.Range("L131").Value = [double]"0.12257617202891925"
$ws.Range("M131").Value = [double]"0.562329950079279"
$ws.Range("L132").Value = [double]"0.72794979301927731"
$ws.Range("M132").Value = [double]"0.12215830376508424"
$ws.Range("L133").Value = [double]"0.4936976297241733"
$ws.Range("M133").Value = [double]"0.47433809487059181"
$ws.Range("L134").Value = [double]"0.082089240811563477"
$ws.Range("M134").Value = [double]"0.52237197208515951"
$ws.Range("L135").Value = [double]"0.40618658330737034"
$ws.Range("M135").Value = [double]"1.3636504794025868"
$ws.Range("L136").Value = [double]"0.3341390449353463"
$ws.Range("M136").Value = [double]"0.94825691301864756"
$ws.Range("L137").Value = [double]"0.17677379446018152"
$ws.Range("M137").Value = [double]"1.6583266532545307"
$ws.Range("L138").Value = [double]"0.36203513686503186"
$ws.Range("M138").Value = [double]"0.84386113340087943"
$ws.Range("L139").Value = [double]"0.34149995374154446"
$ws.Range("M139").Value = [double]"0.70252665323426145"
$ws.Range("L140").Value = [double]"0.6320520016240313"
$ws.Range("M140").Value = [double]"1.1958218312221329"
$ws.Range("L141").Value = [double]"0.48189184084663272"
$ws.Range("M141").Value = [double]"0.50089916453717343"
$ws.Range("L142").Value = [double]"0.80110249618607155"
$ws.Range("M142").Value = [double]"0.064038718615196583"
$ws.Range("L143").Value = [double]"0.060164765283244416"
$ws.Range("M143").Value = [double]"3.6725163303133921"
$ws.Range("L144").Value = [double]"0.77365506304573795"
$ws.Range("M144").Value = [double]"0.89939329340679897"
$ws.Range("L145").Value = [double]"0.33220325716425259"
$ws.Range("M145").Value = [double]"0.95594288841196173"
$ws.Range("L146").Value = [double]"0.20800954417586604"
$ws.Range("M146").Value = [double]"1.6205418802128164"
$ws.Range("L147").Value = [double]"0.90114344332631335"
$ws.Range("M147").Value = [double]"1.0480394104557436"
$ws.Range("L148").Value = [double]"0.063107596731112658"
$ws.Range("M148").Value = [double]"3.5878796546710938"
$ws.Range("L149").Value = [double]"0.01201137181124075"
$ws.Range("M149").Value = [double]"0.38811025082555695"
$ws.Range("L150").Value = [double]"0.94322146206660062"
$ws.Range("M150").Value = [double]"0.0051158929225069288"
$ws.Range("L151").Value = [double]"0.55882358965709955"
$ws.Range("M151").Value = [double]"1.2438987805502757"
$ws.Range("L152").Value = [double]"0.22439810805893068"
$ws.Range("M152").Value = [double]"1.5074923556771453"
$ws.Range("L153").Value = [double]"0.25100802806919165"
$ws.Range("M153").Value = [double]"0.65268791273045812"
$ws.Range("L154").Value = [double]"0.0037576104027197636"
$ws.Range("M154").Value = [double]"0.33352473709515634"
$ws.Range("L155").Value = [double]"0.75648257920294903"
$ws.Range("M155").Value = [double]"0.097061923005201695"

# Set column widths for L and M to match target layout (closest achievable snap)
$ws.Columns.Item(12).ColumnWidth = [double]"14.8"
$ws.Columns.Item(13).ColumnWidth = [double]"14.8"

